# Automatische test-sync: 2025-06-26 21:06:50
# Adds a new log row (15) to the "Logs" sheet for the "MailMind Test"
# mailbox, extends the conditional-formatting ranges to cover it, and
# bumps the "Bestelling / Levering" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new row of log data -----------------------------------
$logs.Range("A15").Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Range("B15").Value = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Range("C15").Value = "Hoi Johan, `nZou je 10 dozen schroeven kunnen bestellen voor de werkplaats? `nWe hebben vooral maat 4x40 nodig.`nGroet, `nRick`nSent using {0}"
$logs.Range("D15").Value = "Bestelling / Levering"
$logs.Range("E15").Value = "Bedankt voor je bericht. Ik neem dit z.s.m. in behandeling."
$logs.Range("F15").Value = "2025-06-26 21:06:46"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Ja"

# Undo the engine's auto row-height pinning caused by the multi-line cell
# (C15) so the row keeps the sheet's default (unpinned) height, matching
# every other data row in the sheet.
$logs.Rows.Item(15).AutoFit()

# --- Extend conditional formatting ranges to include the new row ------
foreach ($fc in $logs.Range("D2:D14").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("D2:D15"))
}
foreach ($fc in $logs.Range("G2:G14").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("G2:G15"))
}
foreach ($fc in $logs.Range("H2:H14").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("H2:H15"))
}
foreach ($fc in $logs.Range("I2:I14").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("I2:I15"))
}

# --- Update the Dashboard summary count --------------------------------
$dashboard.Range("B2").Value = 10
